$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.215.72"
$ws.Range("E2").Value = "  +5.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.053.49"
$ws.Range("E3").Value = "  +5.80%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.69"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.76"
$ws.Range("E6").Value = "  +9.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.708"
$ws.Range("E7").Value = "  +17.70%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.770"
$ws.Range("E9").Value = "  +9.85%  "
$ws.Range("E10").Value = "  +7.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000334"
$ws.Range("E11").Value = "  +5.77%  "
$ws.Range("E12").Value = "  +23.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.11"
$ws.Range("E13").Value = "  +9.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.702.98"
$ws.Range("E14").Value = "  +5.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.030.13"
$ws.Range("E15").Value = "  +4.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.40"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.08"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("E18").Value = "  +3.86%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.170.30"
$ws.Range("E20").Value = "  +5.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.38"
$ws.Range("E21").Value = "  +6.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "100.68"
$ws.Range("E22").Value = "  +17.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.66"
$ws.Range("E23").Value = "  +7.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.10"
$ws.Range("E24").Value = "  +9.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.27"
$ws.Range("E25").Value = "  +9.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.56"
$ws.Range("E26").Value = "  +2.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.09"
$ws.Range("E27").Value = "  +6.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.46"
$ws.Range("E28").Value = "  +6.87%  "
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.39"
$ws.Range("E30").Value = "  +21.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.76"
$ws.Range("E31").Value = "  +5.85%  "
$ws.Range("E32").Value = "  +8.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "677.81"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.82"
$ws.Range("E34").Value = "  +10.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.68"
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.78"
$ws.Range("E36").Value = "  +8.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.444"
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0869"
$ws.Range("E38").Value = "  +6.31%  "
$ws.Range("E39").Value = "  +8.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.46"
$ws.Range("E40").Value = "  +3.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0508"
$ws.Range("E41").Value = "  +8.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.18"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("E45").Value = "  +13.27%  "
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.57"
$ws.Range("E48").Value = "  +14.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.09"
$ws.Range("E49").Value = "  +6.29%  "
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000274"
$ws.Range("E51").Value = "  +5.88%  "
